# Apply weekly update to "Hortaliza, Vega Monumental Concepción - Perejil" sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New date (serial) values for D112:D143 (shifted one pair forward, new date at top)
$newDates = @(
    44694, 44694,   # 112,113
    44433, 44433,   # 114,115
    44203, 44203,   # 116,117
    44308, 44308,   # 118,119
    44665, 44665,   # 120,121
    44383, 44383,   # 122,123
    44272, 44272,   # 124,125
    44237, 44237,   # 126,127
    44330, 44330,   # 128,129
    44187, 44187,   # 130,131
    44194, 44194,   # 132,133
    44365, 44365,   # 134,135
    44327, 44327,   # 136,137
    44358, 44358,   # 138,139
    44217, 44217,   # 140,141
    44460, 44460    # 142,143
)

$startRow = 112
for ($i = 0; $i -lt $newDates.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 4).Value = $newDates[$i]
}

# Append two new rows (144, 145) duplicating the previous content of rows 142/143
# Row 144 (Primera quality)
$ws.Cells.Item(144, 1).Value = 11
$ws.Cells.Item(144, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(144, 3).Value = "Bíobío"
$ws.Cells.Item(144, 4).Value = 44607
$ws.Cells.Item(144, 5).Value = 8
$ws.Cells.Item(144, 6).Value = 100112044
$ws.Cells.Item(144, 7).Value = "Perejil"
$ws.Cells.Item(144, 8).Value = "Sin especificar"
$ws.Cells.Item(144, 9).Value = "Primera"
$ws.Cells.Item(144, 10).Value = 200
$ws.Cells.Item(144, 11).Value = 600
$ws.Cells.Item(144, 12).Value = 700
$ws.Cells.Item(144, 13).Value = 650
$ws.Cells.Item(144, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(144, 15).Value = "Región de Ñuble"
$ws.Cells.Item(144, 16).Value = 650
$ws.Cells.Item(144, 17).Value = 1
$ws.Cells.Item(144, 18).Value = "Hortaliza"

# Row 145 (Segunda quality)
$ws.Cells.Item(145, 1).Value = 11
$ws.Cells.Item(145, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(145, 3).Value = "Bíobío"
$ws.Cells.Item(145, 4).Value = 44607
$ws.Cells.Item(145, 5).Value = 8
$ws.Cells.Item(145, 6).Value = 100112044
$ws.Cells.Item(145, 7).Value = "Perejil"
$ws.Cells.Item(145, 8).Value = "Sin especificar"
$ws.Cells.Item(145, 9).Value = "Segunda"
$ws.Cells.Item(145, 10).Value = 100
$ws.Cells.Item(145, 11).Value = 500
$ws.Cells.Item(145, 12).Value = 500
$ws.Cells.Item(145, 13).Value = 500
$ws.Cells.Item(145, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(145, 15).Value = "Región de Ñuble"
$ws.Cells.Item(145, 16).Value = 500
$ws.Cells.Item(145, 17).Value = 1
$ws.Cells.Item(145, 18).Value = "Hortaliza"

# Copy the date cell style (number format) from row 143's D cell to the new rows' D cells
$ws.Range("D143").Copy() | Out-Null
$ws.Range("D144:D145").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false
